$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 120716.1522112524
$ws.Range("C3").Value = 139051.7984401947
$ws.Range("C6").Value = 122584.2638928925
$ws.Range("C7").Value = 143467.8780872832
$ws.Range("C10").Value = 121897.7867472824
$ws.Range("C11").Value = 141028.3171938551
$ws.Range("C14").Value = 124177.2568180901
$ws.Range("C15").Value = 145824.7654788709
$ws.Range("C18").Value = 142022.830156972
$ws.Range("C19").Value = 170923.5954583264
$ws.Range("C22").Value = 141864.4580373701
$ws.Range("C23").Value = 170794.4238972236
$ws.Range("C26").Value = 141874.069189175
$ws.Range("C27").Value = 170804.409701393
$ws.Range("C30").Value = 141715.9680147173
$ws.Range("C31").Value = 170675.3451249267
$ws.Range("C34").Value = 184772.8045993561
$ws.Range("C35").Value = 200441.6653923719
$ws.Range("C38").Value = 186148.520469406
$ws.Range("C39").Value = 204072.0946239115
$ws.Range("C42").Value = 186036.1053038605
$ws.Range("C43").Value = 202931.8747747873
$ws.Range("C46").Value = 187927.231350612
$ws.Range("C47").Value = 207255.2241745333
$ws.Range("C50").Value = 207211.4476988417
$ws.Range("C51").Value = 236146.7467773402
$ws.Range("C54").Value = 207053.0755731615
$ws.Range("C55").Value = 236017.5756465657
$ws.Range("C58").Value = 207062.6867280056
$ws.Range("C59").Value = 236027.5612394723
$ws.Range("C62").Value = 206904.5855535479
$ws.Range("C63").Value = 235898.4966656377
$ws.Range("C66").Value = 234906.859567781
$ws.Range("C67").Value = 249470.2379004707
$ws.Range("C70").Value = 236064.9238989371
$ws.Range("C71").Value = 252786.2787378968
$ws.Range("C74").Value = 236635.6795396588
$ws.Range("C75").Value = 252773.824457818
$ws.Range("C78").Value = 238447.614645082
$ws.Range("C79").Value = 256959.3831505053
$ws.Range("C82").Value = 260773.8126323012
$ws.Range("C83").Value = 289718.9247643672
$ws.Range("C86").Value = 260615.4405126992
$ws.Range("C87").Value = 289589.7534118955
$ws.Range("C90").Value = 260625.0516645042
$ws.Range("C91").Value = 289599.7390148292
$ws.Range("C94").Value = 260466.9504900465
$ws.Range("C95").Value = 289470.6744309675
$ws.Range("C98").Value = 125819.7983852314
$ws.Range("C99").Value = 152295.34210942
$ws.Range("C102").Value = 131660.3520990623
$ws.Range("C103").Value = 158685.4376367659
$ws.Range("C106").Value = 125819.7983852314
$ws.Range("C107").Value = 152295.34210942
$ws.Range("C110").Value = 131660.3520990623
$ws.Range("C111").Value = 158685.4376367659
$ws.Range("C114").Value = 142591.0276293001
$ws.Range("C115").Value = 171514.5651693879
$ws.Range("C118").Value = 141799.4095141054
$ws.Range("C119").Value = 170675.3451249266
$ws.Range("C122").Value = 142591.0276293001
$ws.Range("C123").Value = 171514.5651693879
$ws.Range("C126").Value = 141799.4095141054
$ws.Range("C127").Value = 170675.3451249266
$ws.Range("C130").Value = 190229.2404402616
$ws.Range("C131").Value = 214131.6591120306
$ws.Range("C134").Value = 195648.0359737862
$ws.Range("C135").Value = 219586.4475331511
$ws.Range("C138").Value = 190229.2404402616
$ws.Range("C139").Value = 214131.6591120306
$ws.Range("C142").Value = 195648.0359737862
$ws.Range("C143").Value = 219586.4475331511
$ws.Range("C146").Value = 207814.1791700112
$ws.Range("C147").Value = 236737.716710099
$ws.Range("C150").Value = 207022.5610548165
$ws.Range("C151").Value = 235898.4966656376
$ws.Range("C154").Value = 207814.1791700112
$ws.Range("C155").Value = 236737.716710099
$ws.Range("C158").Value = 207022.5610548165
$ws.Range("C159").Value = 235898.4966656376
$ws.Range("C162").Value = 240340.4944477063
$ws.Range("C163").Value = 263076.5815397532
$ws.Range("C166").Value = 245577.5501851852
$ws.Range("C167").Value = 268186.0062417226
$ws.Range("C170").Value = 240340.4944477063
$ws.Range("C171").Value = 263076.5815397532
$ws.Range("C174").Value = 245577.5501851852
$ws.Range("C175").Value = 268186.0062417226
$ws.Range("C178").Value = 261386.356935341
$ws.Range("C179").Value = 290309.8944754287
$ws.Range("C182").Value = 260594.7388201463
$ws.Range("C183").Value = 289470.6744309674
$ws.Range("C186").Value = 261386.356935341
$ws.Range("C187").Value = 290309.8944754287
$ws.Range("C190").Value = 260594.7388201463
$ws.Range("C191").Value = 289470.6744309674
